$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Added 64 bit designs: refresh the utilization numbers for row 2 (LUT, LUTRAM, FF, BRAM, DSP)
$ws.Range("B2").Value = 82.06202697753906
$ws.Range("C2").Value = 6.0804595947265625
$ws.Range("D2").Value = 24.877819061279297
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 95.45454406738281

# Column F (DSP) now needs the same width as the other "wide" data columns (~11.71875 chars).
# ColumnWidth is snapped to the nearest whole pixel on the character grid, so pick the
# input value whose rounded result lands on the closest achievable width to 11.71875.
$ws.Range("F1").ColumnWidth = 10.833333333333332
